$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete row 10 (the "Andel heltidsstilling" row), shifting rows 11-30 up.
$ws.Rows.Item(10).Delete()

$ws.Range("D1:D1048576").Validation.Delete()
$u = $ws.Range("D30:D1048576, D1:D28")
foreach ($a in $u.Areas) {
    Write-Output $a.Address()
}
$u.Validation.Add(6, 0, 6, 32)
